# Update gh-pages to output generated at 456a3b4
# Applies numeric updates to the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 555
$ws1.Range("F7").Value = 1728
$ws1.Range("F11").Value = 1832
$ws1.Range("F13").Value = 116
$ws1.Range("F14").Value = 428
$ws1.Range("F22").Value = 788
$ws1.Range("G22").Value = 69
$ws1.Range("F25").Value = 249
$ws1.Range("F26").Value = 263

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 555
$ws4.Range("F7").Value = 1728
$ws4.Range("F12").Value = 1832
$ws4.Range("F14").Value = 116
$ws4.Range("F15").Value = 428
$ws4.Range("F23").Value = 788
$ws4.Range("G23").Value = 69
$ws4.Range("F26").Value = 249
$ws4.Range("F27").Value = 263
